$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("trend")

for ($row = 28; $row -le 51; $row++) {
    $n = $row - 23
    $ws.Cells.Item($row, 21).Formula = "=_input!`$H$n"
    $ws.Cells.Item($row, 29).Formula = "=_input!`$I$n"
}
